# Generate Report for Handback
# Update the handoff/handback generation timestamps on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-09-04 19:10:45"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) / "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-09-04 19:10:40"
$wsZhCn.Range("K2").Value = "2016-09-04 19:10:58"

# de-de sheet: "Correspond Handoff Datetime" (H2) / "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-09-04 19:10:45"
$wsDeDe.Range("K2").Value = "2016-09-04 19:11:10"
